# QuestionpoolHtmlCss.xlsx - "Nop sau phan bien" commit
# - Them duoc quiz cho video bai giang (new quiz question values)
# - Them duoc certificate
#
# Applies the cell-level edits captured in the target diff:
#   H2: 4 -> 3
#   F3: (empty) -> "<h2>"   (new shared string)
#   H3: 4 -> 2
#   H5: 4 -> 2
#   I7: 10 -> 2
#   I8: 10 -> 3
# and leaves the cursor/selection on J8, matching the saved file's
# <selection activeCell="J8" sqref="J8"/>.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3
$ws.Range("F3").Value = "<h2>"
$ws.Range("H3").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("I7").Value = 2
$ws.Range("I8").Value = 3

$ws.Range("J8").Select()
